{"js": "// Load all paragraphs in the body so we can locate the CORE COMPETENCIES\n// bullet paragraphs and the final paragraph of the document by their text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- 1. Collapse the three \"CORE COMPETENCIES\" detail paragraphs into one ---\nconst statIdx = items.findIndex((p) =>\n  p.text.indexOf(\"Statistical Analysis & Machine Learning:\") === 0\n);\nconst bigDataIdx = items.findIndex((p) =>\n  p.text.indexOf(\"Big Data & Data Engineering:\") === 0\n);\nconst dataVizIdx = items.findIndex((p) =>\n  p.text.indexOf(\"Data Visualization & Reporting:\") === 0\n);\n\nif (statIdx === -1 || bigDataIdx === -1 || dataVizIdx === -1) {\n  throw new Error(\"Could not locate CORE COMPETENCIES paragraphs\");\n}\n\n// Replace the first paragraph's text with the new condensed summary line.\nitems[statIdx].insertText(\n  \"Statistical Analysis & Machine Learning \u2022 Big Data & Data Engineering \u2022 Data Visualization & Reporting\",\n  Word.InsertLocation.replace\n);\n\n// Remove the other two (now redundant) detail paragraphs.\nitems[bigDataIdx].delete();\nitems[dataVizIdx].delete();\n\nawait context.sync();\n\n// --- 2. Append a new \"TECHNICAL SKILLS\" section at the end of the body ---\n// Insert the three plain body paragraphs first (while the last paragraph in\n// the document is still a normal/body paragraph), then insert the Heading2\n// title immediately before the first of them. Doing it in this order avoids\n// the new body paragraphs inheriting Heading 2 formatting from the title.\nconst body = context.document.body;\n\nconst p1 = body.insertParagraph(\n  \"STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning; Statistical Computing; A/B Testing; Meta-analytical Techniques\",\n  Word.InsertLocation.end\n);\n\nbody.insertParagraph(\n  \"BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Databases; Data Governance; Streaming Data; Data Pipeline Optimization\",\n  Word.InsertLocation.end\n);\n\nbody.insertParagraph(\n  \"DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Statistical Reporting; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Business Intelligence; Client Presentation\",\n  Word.InsertLocation.end\n);\n\nconst heading = p1.insertParagraph(\"TECHNICAL SKILLS\", Word.InsertLocation.before);\nheading.style = \"Heading 2\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Collapse the three \"CORE COMPETENCIES\" detail paragraphs into one ---\n$statIdx = 0\n$bigDataIdx = 0\n$dataVizIdx = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($statIdx -eq 0 -and $t.StartsWith(\"Statistical Analysis & Machine Learning:\")) {\n        $statIdx = $i\n    } elseif ($bigDataIdx -eq 0 -and $t.StartsWith(\"Big Data & Data Engineering:\")) {\n        $bigDataIdx = $i\n    } elseif ($dataVizIdx -eq 0 -and $t.StartsWith(\"Data Visualization & Reporting:\")) {\n        $dataVizIdx = $i\n    }\n}\n\n# Replace the first paragraph's text with the new condensed summary line\n# (keep the paragraph mark intact so surrounding paragraphs are untouched).\n$p = $d.Paragraphs.Item($statIdx).Range\n$p.Text = \"Statistical Analysis & Machine Learning \" + [char]0x2022 + \" Big Data & Data Engineering \" + [char]0x2022 + \" Data Visualization & Reporting\"\n\n# Delete the other two (now redundant) detail paragraphs, including their\n# paragraph marks. Delete the higher index first so the lower index stays valid.\n$d.Paragraphs.Item($dataVizIdx).Range.Delete()\n$d.Paragraphs.Item($bigDataIdx).Range.Delete()\n\n# --- 2. Append a new \"TECHNICAL SKILLS\" section at the end of the document ---\n# Insert the three plain body paragraphs first (while the document still ends\n# on a normal/body paragraph), then insert the Heading2 title immediately\n# before the first of them. This ordering avoids the new body paragraphs\n# inheriting Heading 2 formatting from the title paragraph.\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\n$last.Range.InsertParagraphAfter()\n$idxFirstBody = $d.Paragraphs.Count\n$d.Paragraphs.Item($idxFirstBody).Range.Text = \"STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning; Statistical Computing; A/B Testing; Meta-analytical Techniques\"\n\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = \"BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Databases; Data Governance; Streaming Data; Data Pipeline Optimization\"\n\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()\n$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = \"DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Statistical Reporting; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Business Intelligence; Client Presentation\"\n\n$d.Paragraphs.Item($idxFirstBody).Range.InsertParagraphBefore()\n$heading = $d.Paragraphs.Item($idxFirstBody)\n$heading.Range.Text = \"TECHNICAL SKILLS\"\n$heading.Style = \"Heading 2\"\n\nWrite-Output \"done\"\n"}
